$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear old selection artifact area / set new values
$ws.Range("A1").Value = "1fa"
$ws.Range("B1").Value = "supplier1"
$ws.Range("C1").Value = 50000
$ws.Range("C1").NumberFormat = "0.00"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "30.08.2016"
$ws.Range("F1").NumberFormat = "mm-dd-yy"

$ws.Range("A2").Value = "1ddd"
$ws.Range("B2").Value = "supplier1"
$ws.Range("C2").Value = 50000
$ws.Range("C2").NumberFormat = "0.00"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.09.2016"
$ws.Range("F2").NumberFormat = "mm-dd-yy"

$ws.Range("E3").Select()
